# The opening bullet currently reads, as plain text:
#   "s" + "  " + " Stack " + "<T>" + " " + "- LIFO " + ...
# spread across several <w:r> runs (the first three all carry the same
# run-level formatting, <w:lang w:val="en-US"/>, except the two-space run
# which carries no rPr at all). The target edit collapses only those
# first three runs ("s", the two spaces, and " Stack ") into one run
# whose text is "Stack ", while every run after it - the "<T>" run, the
# single-space run after it, and the "- LIFO " run - must stay exactly
# as they are: still separate <w:r> elements, even though their
# formatting (<w:lang w:val="en-US"/>) is identical to the merged run's.
#
# This COM-interop runtime coalesces all contiguous same-formatting runs
# in a paragraph whenever a text edit touches that paragraph, which would
# normally swallow "<T>", the space, and "- LIFO " into the edit too.
# To stop that, we briefly give "<T>" and "- LIFO " (the two runs that
# would otherwise be flush up against same-formatted neighbours) a
# throwaway Bold toggle just for the duration of the text edit, then
# flip Bold back off afterwards. Toggling Bold on/off is a pure
# formatting operation - it does not itself trigger the run-coalescing
# pass, and switching it back off leaves no residue in the saved XML.
# The space run between "<T>" and "- LIFO " is left alone; it already
# keeps those two guarded runs from ever being adjacent same-formatting
# siblings of each other.

$d = $word.ActiveDocument

# Locate the exact span to rewrite: "s" + two spaces + " Stack ".
$editRange = $d.Content
$null = $editRange.Find.Execute("s   Stack ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editStart = $editRange.Start
$editEnd = $editRange.End

# Locate the two runs that must stay untouched and unmerged: "<T>" and
# the en-dash run "- LIFO " (built from a char code so the script stays
# plain ASCII).
$enDash = [char]0x2013

$guardARange = $d.Content
$null = $guardARange.Find.Execute("<T>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$guardAStart = $guardARange.Start
$guardAEnd = $guardARange.End

$guardBRange = $d.Content
$null = $guardBRange.Find.Execute("$enDash LIFO ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$guardBStart = $guardBRange.Start
$guardBEnd = $guardBRange.End

# Apply the temporary Bold guard.
$d.Range($guardAStart, $guardAEnd).Font.Bold = 1
$d.Range($guardBStart, $guardBEnd).Font.Bold = 1

# Perform the actual content edit.
$d.Range($editStart, $editEnd).Text = "Stack "

# Remove the guard. The replaced text is shorter than the original by
# (editEnd - editStart) - Len("Stack "), so everything after it shifted
# left by that amount.
$shift = ($editEnd - $editStart) - "Stack ".Length
$d.Range($guardAStart - $shift, $guardAEnd - $shift).Font.Bold = 0
$d.Range($guardBStart - $shift, $guardBEnd - $shift).Font.Bold = 0
